# Adds 43 new incident-log rows (2024-05-13, 14:58:17 .. 16:31:51) to Sheet1,
# extending the used range from A1:G71 to A1:G114 (config test 2 and 3 - v.7.5).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("2024-05-13", "14:58:17", "Palet atascado en la curva", "-", "-", "-", "-"),
    @("2024-05-13", "14:58:19", "No atornilla clips", "-", "-", "-", "-"),
    @("2024-05-13", "14:58:22", "Fallo cámara visión", "-", "-", "-", "-"),
    @("2024-05-13", "14:58:29", "Fallo atornillador", "-", "-", "-", "-"),
    @("2024-05-13", "15:05:16", "Secuencia atornillador", "-", "-", "-", "-"),
    @("2024-05-13", "15:09:22", "Secuencia atornillador", "-", "-", "-", "-"),
    @("2024-05-13", "15:50:56", "No pone tornillo", "-", "-", "-", "-"),
    @("2024-05-13", "15:51:00", "-", "Cámara no detecta foam derecho", "-", "-", "-"),
    @("2024-05-13", "15:51:03", "-", "-", "Power atascado en prensa, cuesta sacar", "-", "-"),
    @("2024-05-13", "15:51:07", "-", "-", "-", "Core enganchado", "-"),
    @("2024-05-13", "15:51:13", "-", "-", "-", "Cover atascado", "-"),
    @("2024-05-13", "15:54:01", "-", "-", "Cámara no detecta Power CP", "-", "-"),
    @("2024-05-13", "15:56:55", "-", "No pone tornillo", "-", "-", "-"),
    @("2024-05-13", "15:56:59", "-", "Fallo en elevador", "-", "-", "-"),
    @("2024-05-13", "15:57:06", "-", "-", "-", "No coloca bien el sealling", "-"),
    @("2024-05-13", "15:57:19", "-", "-", "-", "Power atascado en prensa, cuesta sacar", "-"),
    @("2024-05-13", "15:57:46", "-", "-", "-", "La cámara no detecta Busbar", "-"),
    @("2024-05-13", "16:00:23", "-", "No pone tornillo", "-", "-", "-"),
    @("2024-05-13", "16:00:30", "-", "Fallo tornillo", "-", "-", "-"),
    @("2024-05-13", "16:00:32", "-", "Fallo atornillador", "-", "-", "-"),
    @("2024-05-13", "16:02:12", "-", "AOI no detecta pieza", "-", "-", "-"),
    @("2024-05-13", "16:02:16", "-", "-", "Cámara no detecta Pcb", "-", "-"),
    @("2024-05-13", "16:02:29", "-", "-", "Detección de sealling mal puesto", "-", "-"),
    @("2024-05-13", "16:02:37", "-", "-", "No coloca bien el sealling", "-", "-"),
    @("2024-05-13", "16:03:20", "-", "-", "Cámara no detecta busbar", "-", "-"),
    @("2024-05-13", "16:08:27", "-", "Fallo en paletizador", "-", "-", "-"),
    @("2024-05-13", "16:08:30", "-", "No atornilla clips", "-", "-", "-"),
    @("2024-05-13", "16:08:40", "-", "Etiquetadora", "-", "-", "-"),
    @("2024-05-13", "16:08:46", "-", "-", "Etiquetadora", "-", "-"),
    @("2024-05-13", "16:09:02", "-", "-", "No detecta presencia power CP", "-", "-"),
    @("2024-05-13", "16:09:12", "-", "-", "No detecta presencia power CP", "-", "-"),
    @("2024-05-13", "16:10:16", "-", "-", "Cámara no detecta skeleton", "-", "-"),
    @("2024-05-13", "16:10:27", "-", "-", "-", "Atasco tuerca", "-"),
    @("2024-05-13", "16:11:52", "-", "Etiquetadora", "-", "-", "-"),
    @("2024-05-13", "16:12:00", "-", "Palet atascado en la curva", "-", "-", "-"),
    @("2024-05-13", "16:12:07", "-", "Fallo tornillo", "-", "-", "-"),
    @("2024-05-13", "16:31:07", "-", "No atornilla tapa", "-", "-", "-"),
    @("2024-05-13", "16:31:11", "-", "Ascensor no sube", "-", "-", "-"),
    @("2024-05-13", "16:31:29", "-", "Fallo en elevador", "-", "-", "-"),
    @("2024-05-13", "16:31:34", "-", "-", "No detecta presencia power CP", "-", "-"),
    @("2024-05-13", "16:31:38", "-", "-", "Detección de sealling mal puesto", "-", "-"),
    @("2024-05-13", "16:31:49", "Fallo tolva", "-", "-", "-", "-"),
    @("2024-05-13", "16:31:51", "Fallo cámara visión", "-", "-", "-", "-")
)

$startRow = 72

# Column A holds dates formatted as plain text ("YYYY-MM-DD"), not real Excel
# dates. Force text format first so Excel does not auto-convert those strings
# into date serial numbers while the values are being entered.
$endRow = $startRow + $newRows.Count - 1
$dateRange = $ws.Range("A" + $startRow + ":A" + $endRow)
$dateRange.NumberFormat = "@"

$r = $startRow
foreach ($rowValues in $newRows) {
    $c = 1
    foreach ($val in $rowValues) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}

# Restore the default (unstyled) cell style on column A now that the text is
# committed, so no residual number-format style is left behind on those cells.
$dateRange.Style = "Normal"

